$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45192
# (2023-09-23) to 45202 (2023-10-03) for every data row (rows 2-339).
$ws.Range("C2:C339").Value = 45202
